# Update "想去人数" (want-to-go count) figures across the four sheets to the
# freshly scraped numbers (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$vals1 = @{
    2  = 133
    3  = 336
    4  = 192
    5  = 1213
    6  = 439
    7  = 99
    8  = 161
    14 = 173
    15 = 1436
    16 = 535
    17 = 217
    18 = 327
    20 = 774
    21 = 1129
    23 = 1919
    24 = 2606
    25 = 1371
    27 = 17
    28 = 323
    29 = 393
    30 = 1094
    31 = 794
    32 = 1237
    33 = 145
    35 = 772
    36 = 532
    37 = 633
    38 = 807
    40 = 227
}
foreach ($row in $vals1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $vals1[$row]
}

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$vals2 = @{
    15 = 587
}
foreach ($row in $vals2.Keys) {
    $ws2.Cells.Item($row, 6).Value = $vals2[$row]
}

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$vals3 = @{
    2 = 865
}
foreach ($row in $vals3.Keys) {
    $ws3.Cells.Item($row, 6).Value = $vals3[$row]
}

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$vals4 = @{
    2  = 865
    6  = 133
    7  = 336
    8  = 192
    11 = 1213
    12 = 439
    13 = 99
    14 = 161
    20 = 173
    21 = 1436
    22 = 535
    23 = 217
    24 = 327
    26 = 1129
    27 = 2606
    29 = 1371
    34 = 323
    35 = 393
    36 = 1094
    39 = 794
    40 = 1237
    41 = 772
    42 = 532
    43 = 633
    44 = 807
    48 = 227
}
foreach ($row in $vals4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $vals4[$row]
}
